$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K6").Value = -0.2254024683979639

$ws.Range("J7").Value = -0.1253231084953424
$ws.Range("K7").Value = -0.3352267436446591

$ws.Range("I8").Value = 0.2284633975843539
$ws.Range("J8").Value = 0.01855976243503714

$ws.Range("H9").Value = 0.08028600715190851
$ws.Range("I9").Value = -0.1296176279974082

$ws.Range("G10").Value = -0.07715998185224648
$ws.Range("H10").Value = -0.2870636170015632

$ws.Range("F11").Value = 0.4234994746738243
$ws.Range("G11").Value = 0.2135958395245076

$ws.Range("E12").Value = 0.1431415941383551
$ws.Range("F12").Value = -0.06676204101096155

$ws.Range("D13").Value = 0.3151164519833668
$ws.Range("E13").Value = 0.1052128168340501

$ws.Range("C14").Value = 0.009253912237035311
$ws.Range("D14").Value = -0.2006497229122814

$ws.Range("B15").Value = 0.6215838649243215
$ws.Range("C15").Value = 0.4116802297750048

$ws.Range("B16").Value = -0.2766911554241067
